$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "0.999" or "598.18"
# are not auto-converted to numbers by Excel when re-entered.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '67.907.94'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '3.852.90'
$ws.Range('E3').Value = '  -1.20%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '598.18'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = '166.47'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = '3.852.34'
$ws.Range('E7').Value = '  -1.17%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.525'
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('D11').Value = '6.31'
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('D12').Value = '0.457'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('D14').Value = '36.87'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').Value = '4.496.90'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').Value = '3.828.59'
$ws.Range('E16').Value = '  -2.18%  '
$ws.Range('D17').Value = '67.935.10'
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('D18').Value = '18.15'
$ws.Range('E18').Value = '  +7.02%  '
$ws.Range('D19').Value = '7.41'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('D21').Value = '10.78'
$ws.Range('E21').Value = '  -3.59%  '
$ws.Range('D22').Value = '467.30'
$ws.Range('E22').Value = '  -3.72%  '
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('E24').Value = '  -3.07%  '
$ws.Range('D25').Value = '83.18'
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('D26').Value = '2.22'
$ws.Range('E26').Value = '  -1.10%  '
$ws.Range('D27').Value = '12.17'
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = '10.01'
$ws.Range('E29').Value = '  -0.98%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').Value = '4.000.00'
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('D33').Value = '2.30'
$ws.Range('E33').Value = '  -3.13%  '
$ws.Range('D34').Value = '30.97'
$ws.Range('E34').Value = '  -2.96%  '
$ws.Range('D35').Value = '3.824.37'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('D37').Value = '0.139'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').Value = '5.90'
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('D40').Value = '3.24'
$ws.Range('E40').Value = '  +9.02%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').Value = '0.312'
$ws.Range('E42').Value = '  -2.02%  '
$ws.Range('D43').Value = '425.63'
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '47.17'
$ws.Range('E46').Value = '  -2.62%  '
$ws.Range('D47').Value = '8.54'
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('D48').Value = '143.65'
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('D49').Value = '39.74'
$ws.Range('E49').Value = '  +1.78%  '
$ws.Range('E50').Value = '  +10.67%  '
$ws.Range('E51').Value = '  +0.36%  '
